$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 82
$ws.Range("H2").Value = 88
$ws.Range("E10").Value = 721
$ws.Range("F10").Value = 400
$ws.Range("H10").Value = 495
$ws.Range("E11").Value = 480
$ws.Range("F11").Value = 271
$ws.Range("H11").Value = 336
$ws.Range("E12").Value = 735
$ws.Range("F12").Value = 437
$ws.Range("H12").Value = 523
$ws.Range("E13").Value = 173
$ws.Range("F13").Value = 97
$ws.Range("H13").Value = 131
$ws.Range("E14").Value = 146
$ws.Range("E17").Value = 128
$ws.Range("E22").Value = 201
$ws.Range("F22").Value = 110
$ws.Range("H22").Value = 152
$ws.Range("E23").Value = 232
$ws.Range("F23").Value = 118
$ws.Range("H23").Value = 170
$ws.Range("E24").Value = 277
$ws.Range("F24").Value = 162
$ws.Range("H24").Value = 192
$ws.Range("E25").Value = 339
$ws.Range("F25").Value = 187
$ws.Range("H25").Value = 247
$ws.Range("E26").Value = 213
$ws.Range("F26").Value = 124
$ws.Range("H26").Value = 149
$ws.Range("E27").Value = 394
$ws.Range("F27").Value = 216
$ws.Range("H27").Value = 298
$ws.Range("E28").Value = 233
$ws.Range("F28").Value = 114
$ws.Range("H28").Value = 166
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 118
$ws.Range("H29").Value = 159
$ws.Range("E30").Value = 264
$ws.Range("F30").Value = 162
$ws.Range("G30").Value = 53
$ws.Range("H30").Value = 215
$ws.Range("E31").Value = 84
$ws.Range("G31").Value = 27
$ws.Range("H31").Value = 66
$ws.Range("E32").Value = 223
$ws.Range("F32").Value = 141
$ws.Range("H32").Value = 179
$ws.Range("E33").Value = 343
$ws.Range("F33").Value = 182
$ws.Range("H33").Value = 273
$ws.Range("E34").Value = 261
$ws.Range("F34").Value = 182
$ws.Range("H34").Value = 220
$ws.Range("E36").Value = 91
$ws.Range("E37").Value = 199
$ws.Range("E39").Value = 208
$ws.Range("F39").Value = 106
$ws.Range("H39").Value = 157
$ws.Range("E40").Value = 315
$ws.Range("F40").Value = 161
$ws.Range("H40").Value = 241
$ws.Range("E41").Value = 449
$ws.Range("F41").Value = 226
$ws.Range("H41").Value = 318
$ws.Range("E42").Value = 474
$ws.Range("F42").Value = 268
$ws.Range("H42").Value = 329
$ws.Range("E43").Value = 147
$ws.Range("F43").Value = 85
$ws.Range("H43").Value = 112
$ws.Range("E44").Value = 385
$ws.Range("F44").Value = 202
$ws.Range("H44").Value = 270
$ws.Range("E45").Value = 186
$ws.Range("F45").Value = 102
$ws.Range("H45").Value = 141
$ws.Range("E46").Value = 399
$ws.Range("E47").Value = 551
$ws.Range("F47").Value = 309
$ws.Range("H47").Value = 401
$ws.Range("E48").Value = 281
$ws.Range("E49").Value = 343
$ws.Range("F49").Value = 169
$ws.Range("H49").Value = 256
$ws.Range("E50").Value = 288
$ws.Range("F50").Value = 159
$ws.Range("H50").Value = 232
